$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report record was inserted as row 182 ("Región del
# Maule", 2021-12-29). All the existing records that were at row 182
# onward shift down by one row (to 183..244).
$ws.Rows(182).Insert()

$ws.Range("A182").Value = 10
$ws.Range("B182").Value = "Vega Modelo de Temuco"
$ws.Range("C182").Value = "La Araucanía"
$ws.Range("D182").Value = 44559
$ws.Range("E182").Value = 9
$ws.Range("F182").Value = 100114013
$ws.Range("G182").Value = "Zanahoria"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 220
$ws.Range("K182").Value = 7000
$ws.Range("L182").Value = 8000
$ws.Range("M182").Value = 7455
$ws.Range("N182").Value = "`$/saco 20 kilos"
$ws.Range("O182").Value = "Región del Maule"
$ws.Range("P182").Value = 373
$ws.Range("Q182").Value = 20
$ws.Range("R182").Value = "Hortaliza"
